$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8102118372917175
$ws.Range("B1").Value = 1.016561627388
$ws.Range("C1").Value = 1.468266248703003
$ws.Range("D1").Value = 4.624202728271484
$ws.Range("E1").Value = 3.954318761825562
